# Add the season record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled the same as the rest of the header row (copy the
# format from the last existing header cell so the new ones share style s="1").
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player (rows 2-52) shares the team's overall season record.
$ws.Range("AD2:AD52").Value = 101
$ws.Range("AE2:AE52").Value = 61
$ws.Range("AF2:AF52").Value = 0
